# Apply the "Refined metadata to be additional tab" change:
#  1. Update the panel_query_time-derived "time_taken" timestamps on the
#     existing "data" sheet (F2:F7) to the later re-query time.
#  2. Add a new "metadata" worksheet (after "data") describing the panel
#     query itself (data_name/data_id/data_version/.../panel_get_request).

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# ---------------------------------------------------------------------------
# 1. Refresh the "time_taken" column on the "data" sheet.
# ---------------------------------------------------------------------------
$dataSheet.Cells.Item(2, 6).Value = "2021-10-05 14:34:35.599649"
$dataSheet.Cells.Item(3, 6).Value = "2021-10-05 14:34:35.599657"
$dataSheet.Cells.Item(4, 6).Value = "2021-10-05 14:34:35.599660"
$dataSheet.Cells.Item(5, 6).Value = "2021-10-05 14:34:35.599663"
$dataSheet.Cells.Item(6, 6).Value = "2021-10-05 14:34:35.599666"
$dataSheet.Cells.Item(7, 6).Value = "2021-10-05 14:34:35.599669"

# ---------------------------------------------------------------------------
# 2. Add the "metadata" sheet right after "data".
# ---------------------------------------------------------------------------
$metaSheet = $wb.Worksheets.Add($null, $dataSheet)
$metaSheet.Name = "metadata"

# Header row (B1:G1) re-uses the same bold/bordered/centered style as the
# "data" sheet's own header row -- copy the formatting over instead of
# re-building it by hand so no *new* style entries get minted.
$dataSheet.Cells.Item(1, 2).Copy()
$metaSheet.Range("B1:G1").PasteSpecial(-4122)  # xlPasteFormats

$metaSheet.Cells.Item(1, 2).Value = "data_name"
$metaSheet.Cells.Item(1, 3).Value = "data_id"
$metaSheet.Cells.Item(1, 4).Value = "data_version"
$metaSheet.Cells.Item(1, 5).Value = "data_version_created"
$metaSheet.Cells.Item(1, 6).Value = "panel_query_time"
$metaSheet.Cells.Item(1, 7).Value = "panel_get_request"

# A2 picks up the same style as the "data" sheet's own index column (A2).
$dataSheet.Cells.Item(2, 1).Copy()
$metaSheet.Cells.Item(2, 1).PasteSpecial(-4122)  # xlPasteFormats
$metaSheet.Cells.Item(2, 1).Value = 0

$metaSheet.Cells.Item(2, 2).Value = "Medulloblastoma"
$metaSheet.Cells.Item(2, 3).Value = 3280

# data_version ("0.3") must stay text, not be coerced to a number.
$metaSheet.Cells.Item(2, 4).NumberFormat = "@"
$metaSheet.Cells.Item(2, 4).Value = "0.3"

$metaSheet.Cells.Item(2, 5).Value = "2020-08-10T07:04:41.354417Z"
$metaSheet.Cells.Item(2, 6).Value = "2021-10-05 14:34:35.595777"
$metaSheet.Cells.Item(2, 7).Value = "https://panelapp.agha.umccr.org/api/v1/panels/3280/?format=json"

$dataSheet.Select()
